$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column - copy formatting from neighboring header cell (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
